$p = $ppt.ActivePresentation

# The edited content lives on the slide whose SlideID is 261 (sldId="261" in the
# diff's pc:sldMk) -- that is slide 4 in this deck's Slides collection.
$s = $p.Slides.Item(4)

# The affected shape is the table "Table 6" (shape id 25, creationId
# {3A91F5B0-3974-A14D-A146-FB590F2AAD18} per the ac:graphicFrameMk in the diff).
$tbl = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 25 -and $shp.HasTable) {
        $tbl = $shp.Table
        break
    }
}

# Footnote cell (merged across all 4 columns of row 3) holding the two
# sentences that picked up trailing periods.
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange
$paras = $tr.Paragraphs()

# "Language support is only available in English and Japanese " ->
# "...Japanese. " (period added before the trailing space)
$paras.Item(1).Text = "Language support is only available in English and Japanese. "

# " 1 P2, P3, P4 cases are limited to business hours only in Japan" ->
# "...Japan." (period appended at the very end)
$paras.Item(3).Text = " 1 P2, P3, P4 cases are limited to business hours only in Japan."
